# Update image path prefix in column A ("image" column of Table22) from
# "REPSWITCH1_Practice/" to "Pictures_Practice/" for all data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $val.ToString().StartsWith("REPSWITCH1_Practice/")) {
        $cell.Value2 = $val.ToString().Replace("REPSWITCH1_Practice/", "Pictures_Practice/")
    }
}
